# "Turned on all checks." — the OpenBrowserNotUsed check row in the
# 워크플로우 (Workflow) sheet was left disabled ("No") with a stale check
# file path; flip it to enabled ("Yes") and point it at the check file's
# current location.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(17, 1).Value = "Yes"
$ws.Cells.Item(17, 3).Value = "Checks\Standard\UnusedVariables\UnusedVariables.xaml"

# Keep the workflow sheet active/selected at A1 (matches the saved view
# state after the edit).
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
